# Apply roster/team updates to the NSY sheet.
# The underlying rows were reshuffled (players/teams swapped between rows),
# so we overwrite the Player / Position / Team columns for each affected
# row with their new final values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Jalen Williams, SG,SF,PF,C, Oklahoma City Thunder  ->  RJ Barrett, SG,SF,PF, Toronto Raptors
$ws.Range("A7").Value = "RJ Barrett"
$ws.Range("B7").Value = "SG,SF,PF"
$ws.Range("C7").Value = "Toronto Raptors"

# Row 10: John Collins, PF,C, Utah Jazz  ->  Jalen Williams, SG,SF,PF,C, Oklahoma City Thunder
$ws.Range("A10").Value = "Jalen Williams"
$ws.Range("B10").Value = "SG,SF,PF,C"
$ws.Range("C10").Value = "Oklahoma City Thunder"

# Row 12: Jimmy Butler, SF,PF, Miami Heat  ->  Jusuf Nurkic, C, Phoenix Suns
$ws.Range("A12").Value = "Jusuf Nurkic"
$ws.Range("B12").Value = "C"
$ws.Range("C12").Value = "Phoenix Suns"

# Row 13: RJ Barrett, SG,SF,PF, Toronto Raptors  ->  Jimmy Butler, SF,PF, Miami Heat
$ws.Range("A13").Value = "Jimmy Butler"
$ws.Range("B13").Value = "SF,PF"
$ws.Range("C13").Value = "Miami Heat"

# Row 14: Terry Rozier, PG, Miami Heat  ->  John Collins, PF,C, Utah Jazz
$ws.Range("A14").Value = "John Collins"
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "Utah Jazz"

# Row 16: Kevin Huerter, SG,SF, Sacramento Kings  ->  Draymond Green, PF,C, Golden State Warriors
$ws.Range("A16").Value = "Draymond Green"
$ws.Range("B16").Value = "PF,C"
$ws.Range("C16").Value = "Golden State Warriors"
